$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Despesa"
$ws.Range("B13").Value = "SERVIÇOS"
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = "30/01/2025"
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = " "
